$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the label columns (A:C) to fit the new content.
$ws.Columns("A:C").ColumnWidth = 33.45

# Row height tweaks that came with the new column.
$ws.Rows(1).RowHeight = 40.5
$ws.Rows(2).RowHeight = 13.5
$ws.Rows(3).RowHeight = 13.5

# Clear the lingering N6 selection left over in the saved view.
$ws.Range("A1").Select() | Out-Null

# Bring column N (2023) into existence with the same look as column M (2022),
# then overwrite the values with the new year's figures.
$ws.Range("M3:M11").Copy() | Out-Null
$ws.Range("N3:N11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("N4").Value = 2023
$ws.Range("N5").Value = 0.86
$ws.Range("N6").Value = 1.07
$ws.Range("N7").Value = 25.27
$ws.Range("N8").Value = 14
$ws.Range("N9").Value = 0.12
$ws.Range("N10").Value = 21.74
$ws.Range("N11").Value = 9.4600000000000009
